# Apply a straightforward case-fix rename: "cyclones" -> "Cyclones"
# across the sheet tab name and every cell value that held that team name.

$wb = $excel.ActiveWorkbook

# 1) Rename the worksheet tab itself ("cyclones" -> "Cyclones")
$cyclonesSheet = $wb.Worksheets.Item("cyclones")
$cyclonesSheet.Name = "Cyclones"

# 2) Update the "League Summary" sheet: team name cell A2
$leagueSummary = $wb.Worksheets.Item("League Summary")
$leagueSummary.Range("A2").Value = "Cyclones"

# 3) Update the "Player Summary" sheet: team name column B, rows 2-12
$playerSummary = $wb.Worksheets.Item("Player Summary")
$lastRow = $playerSummary.Cells.Item($playerSummary.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $playerSummary.Cells.Item($r, 2)
    if ($cell.Value2 -eq "cyclones") {
        $cell.Value = "Cyclones"
    }
}
